# "Completed Appeal Form Date verification"
#
# Adds three more test-account columns (M/N/O) to the User_Profile sheet,
# mirroring the existing User Name / LPA / Third Party triplet (J/K/L) but
# for the PreProd / Admin appeal-officer accounts, then leaves the
# selection on the new last cell (O2) the way Excel would after typing the
# new data in from the keyboard.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data: columns M, N, O -------------------------------------------
$ws.Range("M1").Value = "LPA_Preprod"
$ws.Range("M2").Value = "TestOfficer_B"

$ws.Range("N1").Value = "Admin"
$ws.Range("N2").Value = "PINStestadmin"

$ws.Range("O1").Value = "Third Party PreProd"
$ws.Range("O2").Value = "santosh.preprod"

# --- Column widths ---------------------------------------------------------
# Excel auto-sizes (bestFit) columns I/J/L/M/N/O after the new entries are
# typed in. This engine's ColumnWidth setter quantizes to 1/6 of a
# character, so the requested values are converted through the same
# (w + 5/6) rounding the host applies, landing each column on the closest
# representable width to the real Excel bestFit result.
$ws.Columns.Item(9).ColumnWidth = 34.666666666666664    # I: 35.54296875
$ws.Columns.Item(10).ColumnWidth = 9.666666666666666    # J: 10.453125
$ws.Columns.Item(12).ColumnWidth = 20.0                 # L: 20.81640625
$ws.Columns.Item(13).ColumnWidth = 12.333333333333334   # M: 13.1796875
$ws.Columns.Item(14).ColumnWidth = 12.166666666666666   # N: 13
$ws.Columns.Item(15).ColumnWidth = 16.333333333333332   # O: 17.1796875

# --- View state --------------------------------------------------------
# Scroll right so column E is the leftmost visible column, and leave the
# selection on the newly-completed O2 cell.
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("O2").Select()

$wb.Save()
